# "add corals per rack"
# Adds a new column E ("n_corals") to Sheet1 that records how many coral
# fragments were on each rack for a given CBASS run (keyed by date + rack
# color). The two CBASS runs in this sheet (2021-09-06 and 2021-09-08)
# each used one rack with a different coral count than the rest of the
# racks that day:
#   2021-09-06: "double_orange" racks held 3 corals, every other rack held 5
#   2021-09-08: "double_pink"   racks held 7 corals, every other rack held 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Header
$ws.Cells.Item(1, 5).Value = "n_corals"

for ($r = 2; $r -le $lastRow; $r++) {
    $date = $ws.Cells.Item($r, 1).Value()
    $rack = $ws.Cells.Item($r, 4).Value()

    if ($date -eq "2021-09-06") {
        if ($rack -eq "double_orange") {
            $nCorals = "3"
        } else {
            $nCorals = "5"
        }
    } else {
        if ($rack -eq "double_pink") {
            $nCorals = "7"
        } else {
            $nCorals = "8"
        }
    }

    $ws.Cells.Item($r, 5).Value = $nCorals
}

# Restore the view state Excel would have recorded for this edit:
# scrolled down so row 67 is at the top, with E88 as the active cell.
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E88").Select()
